$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlinks from the "UserName" column (B2:B5) before
# the columns are deleted so they don't reattach to the shifted data.
$ws.Cells.Hyperlinks.Delete()

# Remove the "UserName" and "Password" columns (B:C); the remaining
# columns (NewResourceTitle..TextToAdd) shift left to B:G.
$ws.Columns("B:C").Delete()

# Restore the active selection used when the sheet was last saved.
$ws.Range("D12").Select()
